# Applies the daily refresh of the Terminal Status Cassette Balances report:
# updated cash-balance figures/timestamps for 11/08/23, plus three rows of
# terminals that dropped out of (or re-entered) the "needs attention" list,
# which shifted the row order for terminal IDs L662336 .. L682801.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: SCL ENTERPRISES LAUNDRY (LK644532) ---
$ws.Range("E6").Value = 2200
$ws.Range("J6").Value = "11/08/23 09:37"
$ws.Range("K6").Value = "11/08/23 09:37"
$ws.Range("M6").Value = "`$2,200 as of 11/8/2023 7:37:51 AM"
$ws.Range("N6").Value = 2240

# --- Row 9: NICK SHELL SERVICE (L474792) ---
$ws.Range("E9").Value = 2540
$ws.Range("J9").Value = "11/08/23 12:46"
$ws.Range("K9").Value = "11/08/23 12:46"
$ws.Range("M9").Value = "`$2,540 as of 11/8/2023 10:46:48 AM"
$ws.Range("N9").Value = 2600

# --- Row 10: now L662336 / SB#4 MONA MARKET (was LK561655 / CRENSHAW CRAVOR #2) ---
$ws.Range("A10").Value = "L662336"
$ws.Range("C10").Value = "SB#4 MONA MARKET"
$ws.Range("E10").Value = 2560
$ws.Range("I10").ClearContents()
$ws.Range("H10").Value = 45250.0418602662
$ws.Range("J10").Value = "11/08/23 17:02"
$ws.Range("K10").Value = "11/08/23 17:02"
$ws.Range("L10").Value = 120
$ws.Range("M10").Value = "`$3,180 as of 11/8/2023 7:54:25 AM"
$ws.Range("N10").Value = 2640

# --- Row 11: now LK561655 / CRENSHAW CRAVOR #2 (was L474746 / ZACATES MARKET) ---
$ws.Range("A11").Value = "LK561655"
$ws.Range("C11").Value = "CRENSHAW CRAVOR #2"
$ws.Range("E11").Value = 2780
$ws.Range("H11").ClearContents()
$ws.Range("I11").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J11").Value = "01/23/20 08:24"
$ws.Range("K11").Value = "01/23/20 08:24"
$ws.Range("M11").Value = "`$2,780 as of 1/23/2020 6:24:32 AM"
$ws.Range("N11").Value = 2800

# --- Row 12: now L475090 / S.B. 2 (was L662336 / SB#4 MONA MARKET) ---
$ws.Range("A12").Value = "L475090"
$ws.Range("C12").Value = "S.B. 2"
$ws.Range("E12").Value = 2840
$ws.Range("H12").Value = 45241.0418602662
$ws.Range("J12").Value = "11/08/23 19:47"
$ws.Range("K12").Value = "11/08/23 19:03"
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = "`$3,220 as of 11/8/2023 9:27:02 AM"
$ws.Range("N12").Value = 2840

# --- Row 13: now L474746 / ZACATES MARKET (was L475090 / S.B. 2) ---
$ws.Range("A13").Value = "L474746"
$ws.Range("C13").Value = "ZACATES MARKET"
$ws.Range("E13").Value = 2960
$ws.Range("H13").Value = 45247.0418602662
$ws.Range("J13").Value = "11/08/23 13:39"
$ws.Range("K13").Value = "11/08/23 13:36"
$ws.Range("M13").Value = "`$2,960 as of 11/8/2023 11:36:06 AM"
$ws.Range("N13").Value = 2960

# --- Row 14: now LK864765 / SKY LIQUOR (was L475182 / LA ESQUINA DE ORO) ---
$ws.Range("A14").Value = "LK864765"
$ws.Range("C14").Value = "SKY LIQUOR"
$ws.Range("E14").Value = 3400
$ws.Range("I14").ClearContents()
$ws.Range("H14").Value = 45243.0418602662
$ws.Range("J14").Value = "11/08/23 17:12"
$ws.Range("K14").Value = "11/08/23 17:12"
$ws.Range("L14").Value = 80
$ws.Range("M14").Value = "`$3,840 as of 11/8/2023 11:37:39 AM"
$ws.Range("N14").Value = 3500

# --- Row 15: now L475182 / LA ESQUINA DE ORO (was L688961 / MONA MART) ---
$ws.Range("A15").Value = "L475182"
$ws.Range("C15").Value = "LA ESQUINA DE ORO"
$ws.Range("E15").Value = 3800
$ws.Range("I15").Value = "ATM Inactive greater than 48 minutes"
$ws.Range("J15").Value = "09/16/20 16:57"
$ws.Range("K15").Value = "09/15/20 23:38"
$ws.Range("M15").Value = "`$3,800 as of 9/16/2020 1:28:00 PM"
$ws.Range("N15").Value = 3800

# --- Row 16: S B MARKET ST (L697590) ---
$ws.Range("E16").Value = 3820
$ws.Range("J16").Value = "11/08/23 19:38"
$ws.Range("K16").Value = "11/08/23 17:11"
$ws.Range("M16").Value = "`$4,020 as of 11/8/2023 10:06:15 AM"
$ws.Range("N16").Value = 4020

# --- Row 17: now L688966 / S B WESTERN 108TH MARKET (was LK864765 / SKY LIQUOR) ---
$ws.Range("A17").Value = "L688966"
$ws.Range("C17").Value = "S B WESTERN 108TH MARKET"
$ws.Range("E17").Value = 3880
$ws.Range("H17").Value = 45245.0418602662
$ws.Range("J17").Value = "11/08/23 19:10"
$ws.Range("K17").Value = "11/08/23 19:10"
$ws.Range("M17").Value = "`$4,360 as of 11/8/2023 10:31:03 AM"
$ws.Range("N17").Value = 3840

# --- Row 18: now L688961 / MONA MART (was L688966 / S B WESTERN 108TH MARKET) ---
$ws.Range("A18").Value = "L688961"
$ws.Range("C18").Value = "MONA MART"
$ws.Range("E18").Value = 4000
$ws.Range("H18").ClearContents()
$ws.Range("I18").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J18").Value = "10/17/23 13:26"
$ws.Range("K18").Value = "10/17/23 13:00"
$ws.Range("L18").Value = 0
$ws.Range("M18").Value = "`$4,000 as of 10/17/2023 11:00:09 AM"
$ws.Range("N18").Value = 4000

# --- Row 19: WORLDWIDE AUTOMOTIVE (LK236828) ---
$ws.Range("E19").Value = 4740
$ws.Range("J19").Value = "11/08/23 17:48"
$ws.Range("K19").Value = "11/08/23 17:48"
$ws.Range("M19").Value = "`$4,840 as of 11/8/2023 10:15:56 AM"
$ws.Range("N19").Value = 4840

# --- Row 20: BABS MARKET (L474761) ---
$ws.Range("E20").Value = 4960
$ws.Range("J20").Value = "11/08/23 17:33"
$ws.Range("K20").Value = "11/08/23 17:33"
$ws.Range("N20").Value = 5120

# --- Row 21: SAFETY MARKET (L474817) ---
$ws.Range("E21").Value = 5660
$ws.Range("J21").Value = "11/08/23 18:26"
$ws.Range("K21").Value = "11/08/23 18:26"
$ws.Range("M21").Value = "`$5,820 as of 11/8/2023 10:55:36 AM"
$ws.Range("N21").Value = 5760

# --- Row 22: now L704741 / W ADAMS COIN LAUNDRY (was L682801 / SB#5) ---
$ws.Range("A22").Value = "L704741"
$ws.Range("C22").Value = "W ADAMS COIN LAUNDRY"
$ws.Range("E22").Value = 7780
$ws.Range("I22").ClearContents()
$ws.Range("H22").Value = 45248.0418602662
$ws.Range("J22").Value = "11/08/23 17:29"
$ws.Range("K22").Value = "11/08/23 17:29"
$ws.Range("L22").Value = 40
$ws.Range("M22").Value = "`$8,040 as of 11/8/2023 11:37:41 AM"
$ws.Range("N22").Value = 7800

# --- Row 23: DONUT & SANDWICH (L476340) ---
$ws.Range("E23").Value = 7840
$ws.Range("J23").Value = "11/08/23 11:55"
$ws.Range("K23").Value = "11/08/23 11:55"
$ws.Range("M23").Value = "`$7,840 as of 11/8/2023 9:55:43 AM"
$ws.Range("N23").Value = 7860

# --- Row 24: now L682801 / SB#5 (was L704741 / W ADAMS COIN LAUNDRY) ---
$ws.Range("A24").Value = "L682801"
$ws.Range("C24").Value = "SB#5"
$ws.Range("E24").Value = 7840
$ws.Range("H24").ClearContents()
$ws.Range("I24").Value = "ATM Inactive greater than 2000 minutes"
$ws.Range("J24").Value = "09/28/23 15:22"
$ws.Range("K24").Value = "09/28/23 12:14"
$ws.Range("L24").Value = 0
$ws.Range("M24").Value = "`$7,840 as of 9/28/2023 12:31:50 PM"
$ws.Range("N24").Value = 7840

# --- Row 25: SAMYS PHONE CARDS (LK923383) ---
$ws.Range("E25").Value = 10460
$ws.Range("J25").Value = "11/08/23 17:39"
$ws.Range("K25").Value = "11/08/23 17:39"
$ws.Range("M25").Value = "`$10,520 as of 11/8/2023 11:38:10 AM"
$ws.Range("N25").Value = 10500

# --- Row 26: 98 DISCOUNT STORE (LK891176) ---
$ws.Range("E26").Value = 10920
$ws.Range("J26").Value = "11/08/23 19:37"
$ws.Range("K26").Value = "11/08/23 19:14"
$ws.Range("M26").Value = "`$11,360 as of 11/8/2023 11:41:50 AM"
$ws.Range("N26").Value = 10920

# --- Row 27: S B DISCOUNT MART (L697589) ---
$ws.Range("E27").Value = 13920
$ws.Range("J27").Value = "11/08/23 19:03"
$ws.Range("K27").Value = "11/08/23 19:03"
$ws.Range("M27").Value = "`$14,140 as of 11/8/2023 11:40:47 AM"
$ws.Range("N27").Value = 14120

# --- Row 28: Total Outstanding Cash Balance ---
$ws.Range("E28").Value = 115540
